$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs that were split by the "_GoBack" bookmark in the
#    "Worked on Application to Container (A2C) migration..." paragraph back
#    into a single run / single piece of text. A whole-phrase Find & Replace
#    (old text == new text) naturally collapses the split runs into one run
#    and removes the bookmark that used to sit between them.
# ---------------------------------------------------------------------------
$oldCombined = "Worked on Application to Container (A2C) migration. It involves migrating applications running on virtual machine to Kubernetes cluster. Technology used: Python, Docker, Kubernetes, SaaS"
$rngMerge = $d.Content
$rngMerge.Find.Execute($oldCombined, $true, $false, $false, $false, $false, $true, 1, $false, $oldCombined, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Swap the text of the two "Link: ..." project paragraphs:
#      - the first one ("electricity_bill_collection") becomes "Question_Comparator"
#      - the second one ("Question_Comparator") becomes "electricity_bill_collection"
#    Both substrings are unique in the document, so a plain Find & Replace
#    (via a temporary placeholder, to avoid the second replace re-matching
#    text the first replace just produced) is safe.
# ---------------------------------------------------------------------------
$rngA = $d.Content
$rngA.Find.Execute("https://github.com/aayushsinha44/electricity_bill_collection", $true, $false, $false, $false, $false, $true, 1, $false, "https://github.com/aayushsinha44/___TMP_SWAP___", 2) | Out-Null

$rngB = $d.Content
$rngB.Find.Execute("https://github.com/aayushsinha44/Question_Comparator", $true, $false, $false, $false, $false, $true, 1, $false, "https://github.com/aayushsinha44/electricity_bill_collection", 2) | Out-Null

$rngC = $d.Content
$rngC.Find.Execute("https://github.com/aayushsinha44/___TMP_SWAP___", $true, $false, $false, $false, $false, $true, 1, $false, "https://github.com/aayushsinha44/Question_Comparator", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark as a collapsed point right after the
#    text of the (now first) "Link: .../Question_Comparator" paragraph -
#    i.e. after the run, before the paragraph mark - matching the target.
#
#    Range.Collapse()/SetRange()/MoveStart() produce a stale anchor that
#    Bookmarks.Add silently mis-places, so instead: temporarily append a
#    unique marker run right after the target text, wrap a fresh Find hit on
#    that marker with the bookmark, then delete the marker text again. Since
#    the whole bookmarked span is removed, Word collapses the bookmark to a
#    single point exactly where the marker used to be.
# ---------------------------------------------------------------------------
$marker = "ZZ_GOBACK_MARKER_ZZ"

$rngFind = $d.Content
$rngFind.Find.Execute("Link: https://github.com/aayushsinha44/Question_Comparator", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$rngFind.InsertAfter($marker)

$rngMarker = $d.Content
$rngMarker.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $rngMarker) | Out-Null

$rngCleanup = $d.Content
$rngCleanup.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
